$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2708946666666667
$ws.Range("H2").Value = 0.812684
$ws.Range("I2").Value = 0.1616296696421007
$ws.Range("J2").Value = 0.1616296696421007
$ws.Range("M2").Value = 0.3045636666666667
$ws.Range("N2").Value = 0.913691
$ws.Range("Q2").Value = 0.08250467296044445
$ws.Range("R2").Value = 0.742542056644
$ws.Range("S2").Value = 0.1616296696421007
$ws.Range("T2").Value = 0.1616296696421007

# Row 3
$ws.Range("I3").Value = 0.6313295261673385
$ws.Range("J3").Value = 0.6313295261673384
$ws.Range("M3").Value = 0.3045636666666667
$ws.Range("N3").Value = 0.913691
$ws.Range("Q3").Value = 0.3222653130582223
$ws.Range("R3").Value = 2.900387817524
$ws.Range("S3").Value = 0.6313295261673385
$ws.Range("T3").Value = 0.6313295261673384

# Row 4
$ws.Range("I4").Value = 0.2070408041905609
$ws.Range("J4").Value = 0.2070408041905609
$ws.Range("M4").Value = 0.3045636666666667
$ws.Range("N4").Value = 0.913691
$ws.Range("Q4").Value = 0.1056850136304445
$ws.Range("R4").Value = 0.9511651226740001
$ws.Range("S4").Value = 0.2070408041905609
$ws.Range("T4").Value = 0.2070408041905609
